$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.2434456928838951
$ws.Range("C2").Value = 0.4868913857677903
$ws.Range("J2").Value = 0.01872659176029963
$ws.Range("P2").Value = 0.1685393258426966
$ws.Range("S2").Value = 0.08239700374531835
$ws.Range("B3").Value = 0.01550387596899225
$ws.Range("C3").Value = 0.02325581395348837
$ws.Range("J3").Value = 0.02325581395348837
$ws.Range("P3").Value = 0.7054263565891473
$ws.Range("S3").Value = 0.2325581395348837
$ws.Range("J4").Value = 0.0625
$ws.Range("P4").Value = 0.625
$ws.Range("S4").Value = 0.3125
$ws.Range("B6").Value = 0.08021390374331551
$ws.Range("F6").Value = 0.06951871657754011
$ws.Range("J6").Value = 0.2085561497326203
$ws.Range("O6").Value = 0.0213903743315508
$ws.Range("R6").Value = 0.09625668449197861
$ws.Range("S6").Value = 0.4064171122994653
$ws.Range("B7").Value = 0.07881773399014778
$ws.Range("D7").Value = 0.01970443349753695
$ws.Range("F7").Value = 0.03940886699507389
$ws.Range("J7").Value = 0.1527093596059113
$ws.Range("Q7").Value = 0.1773399014778325
$ws.Range("R7").Value = 0.08866995073891626
$ws.Range("S7").Value = 0.4433497536945813
$ws.Range("B8").Value = 0.07630522088353414
$ws.Range("D8").Value = 0.01807228915662651
$ws.Range("F8").Value = 0.06224899598393574
$ws.Range("J8").Value = 0.08835341365461848
$ws.Range("O8").Value = 0.004016064257028112
$ws.Range("Q8").Value = 0.1867469879518072
$ws.Range("R8").Value = 0.0783132530120482
$ws.Range("S8").Value = 0.4859437751004016
$ws.Range("B9").Value = 0.07865168539325842
$ws.Range("D9").Value = 0.01685393258426966
$ws.Range("F9").Value = 0.0449438202247191
$ws.Range("J9").Value = 0.1067415730337079
$ws.Range("O9").Value = 0.01685393258426966
$ws.Range("Q9").Value = 0.2134831460674157
$ws.Range("R9").Value = 0.1292134831460674
$ws.Range("S9").Value = 0.3932584269662922
$ws.Range("B10").Value = 0.09289176090468497
$ws.Range("D10").Value = 0.01534733441033926
$ws.Range("E10").Value = 0.0008077544426494346
$ws.Range("F10").Value = 0.05654281098546042
$ws.Range("J10").Value = 0.1203554119547657
$ws.Range("O10").Value = 0.005654281098546042
$ws.Range("Q10").Value = 0.2205169628432956
$ws.Range("R10").Value = 0.07350565428109855
$ws.Range("S10").Value = 0.4143780290791599
$ws.Range("G11").Value = 0.1615853658536585
$ws.Range("J11").Value = 0.0975609756097561
$ws.Range("K11").Value = 0.2439024390243902
$ws.Range("L11").Value = 0.4908536585365854
$ws.Range("S11").Value = 0.006097560975609756
$ws.Range("G12").Value = 0.7839506172839507
$ws.Range("J12").Value = 0.1851851851851852
$ws.Range("K12").Value = 0.01234567901234568
$ws.Range("S12").Value = 0.01851851851851852
$ws.Range("G13").Value = 0.8571428571428571
$ws.Range("J13").Value = 0.1142857142857143
$ws.Range("S13").Value = 0.02857142857142857
$ws.Range("F15").Value = 0.0101010101010101
$ws.Range("H15").Value = 0.2121212121212121
$ws.Range("I15").Value = 0.06565656565656566
$ws.Range("J15").Value = 0.4141414141414141
$ws.Range("K15").Value = 0.0707070707070707
$ws.Range("M15").Value = 0.005050505050505051
$ws.Range("O15").Value = 0.04545454545454546
$ws.Range("S15").Value = 0.1767676767676768
$ws.Range("F16").Value = 0.0261437908496732
$ws.Range("H16").Value = 0.2026143790849673
$ws.Range("I16").Value = 0.08496732026143791
$ws.Range("J16").Value = 0.392156862745098
$ws.Range("K16").Value = 0.0718954248366013
$ws.Range("M16").Value = 0.0196078431372549
$ws.Range("N16").Value = 0.006535947712418301
$ws.Range("O16").Value = 0.08496732026143791
$ws.Range("S16").Value = 0.1111111111111111
$ws.Range("F17").Value = 0.008676789587852495
$ws.Range("H17").Value = 0.193058568329718
$ws.Range("I17").Value = 0.08459869848156182
$ws.Range("J17").Value = 0.4273318872017354
$ws.Range("K17").Value = 0.09544468546637744
$ws.Range("M17").Value = 0.008676789587852495
$ws.Range("O17").Value = 0.07592190889370933
$ws.Range("S17").Value = 0.1062906724511931
$ws.Range("F18").Value = 0.03141361256544502
$ws.Range("H18").Value = 0.1570680628272251
$ws.Range("I18").Value = 0.08376963350785341
$ws.Range("J18").Value = 0.4136125654450262
$ws.Range("K18").Value = 0.08900523560209424
$ws.Range("M18").Value = 0.02617801047120419
$ws.Range("O18").Value = 0.06282722513089005
$ws.Range("S18").Value = 0.1361256544502618
$ws.Range("F19").Value = 0.01370906321401371
$ws.Range("H19").Value = 0.2368621477532369
$ws.Range("I19").Value = 0.0753998476770754
$ws.Range("J19").Value = 0.3655750190403656
$ws.Range("K19").Value = 0.1172886519421173
$ws.Range("M19").Value = 0.01675552170601675
$ws.Range("O19").Value = 0.06549885757806551
$ws.Range("S19").Value = 0.1089108910891089
